# Update the public epexspot_prices.xlsx workbook with the latest data:
#  - "Prix Spot" sheet: add a new day column AD (13-jul) with its hourly prices
#  - "Gaz" sheet: add a new row (2025-07-11) with its last price
#  - "CO2" sheet: add a new row (2025-07-11) with its last price

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append column AD ("13-jul") with its 24 hourly prices
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# New header cell AD1 gets the exact same formatting as the other day headers
$wsPrix.Range("AC1").Copy()
$wsPrix.Range("AD1").PasteSpecial(-4122)  # xlPasteFormats
$wsPrix.Range("AD1").Value = "13-jul"

$prixValues = @(95.66,87.25,74.89,60.72,55.25,54.07,51.09,62.31,55,22.5,4.28,0.05,0,0,0,0,1,9.369999999999999,30.39,70.09999999999999,107.69,122.8,126.94,114.4)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 30).Value = $prixValues[$i]
}

$wsPrix.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 27 (2025-07-11)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date-looking text to stay plain text (not auto-converted to a
# date serial number), then restore the plain/default style used by the
# rest of the column A data cells.
$wsGaz.Range("A27").NumberFormat = "@"
$wsGaz.Range("A27").Value = "2025-07-11"
$wsGaz.Range("A27").Style = $wsGaz.Range("A26").Style
$wsGaz.Range("B27").Value = 34.8

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 27 (2025-07-11)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A27").NumberFormat = "@"
$wsCo2.Range("A27").Value = "2025-07-11"
$wsCo2.Range("A27").Style = $wsCo2.Range("A26").Style
$wsCo2.Range("B27").Value = 69.8
